# Huge v2 translation patch update.
#
# The sheet originally stored the source text in column A and, for rows
# that had a translation, the translated text in column D (occasionally
# duplicated into column C as well). The new layout consolidates every
# translation into column B: rows that already had a column D value move
# that value into B, while rows that had no translation yet get column A's
# text copied into B as a placeholder. Columns C and D are then removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 31

for ($r = 1; $r -le $lastRow; $r++) {
    $sourceVal = $ws.Cells.Item($r, 1).Value()
    $dVal = $ws.Cells.Item($r, 4).Value()

    if ($dVal -ne $null -and $dVal -ne "") {
        $ws.Cells.Item($r, 2).Value = $dVal
    } else {
        $ws.Cells.Item($r, 2).Value = $sourceVal
    }
}

# Drop the now-obsolete columns C and D entirely.
$ws.Range("C1:D" + $lastRow).Clear()

# Re-fit row heights so that writing the (often multi-line) translation
# strings doesn't leave stray custom row heights behind.
$ws.Range("A1:B" + $lastRow).Rows.AutoFit()
